$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Near the end of the document, insert a new paragraph (right before
#    the final "Create a feature image..." paragraph) that will carry
#    the bold title text "Play Ace Ventura for Free - Slot Game Review".
#    We build it by splitting off a fresh paragraph after the last
#    "ListBullet" item, resetting it to the Normal style, and then
#    copying the FormattedText of the (still present, top-of-doc) Meta
#    description paragraph into it, so it picks up the same leading
#    empty run / bold run shape -- finally swapping the wording via a
#    scoped Find/Replace.
# ---------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)

$count = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs($count - 1)
$secondToLast.Range.InsertParagraphAfter()

$count = $d.Paragraphs.Count
$newPara = $d.Paragraphs($count - 1)
$newPara.Style = $d.Styles("Normal")

$newRange = $newPara.Range
$newRange.FormattedText = $metaPara.Range.FormattedText
$newRange.Find.Execute("Meta description: Discover the amusing online slot game Ace Ventura and its various bonus features. Play for free and enjoy the immersive experience.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Play Ace Ventura for Free - Slot Game Review", 2)

# ---------------------------------------------------------------------
# 2) Remove the original "Meta description: ..." paragraph that sits
#    right under the H1 title at the top of the document.
# ---------------------------------------------------------------------
$d.Paragraphs(2).Range.Delete()

# ---------------------------------------------------------------------
# 3) Replace the text of the final paragraph (the old image prompt)
#    with the meta-description wording, keeping its italic run format.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastRange = $lastPara.Range
$lastRange.Find.Execute("Create a feature image for Ace Ventura that showcases the game's cartoon style and features a happy Maya warrior with glasses. The image should be colorful and eye-catching, with the Maya warrior standing in the foreground with a big smile, wearing traditional warrior headdress, and holding a magnifying glass in hand. The background should feature elements from the game, such as Ace Ventura characters, animals, or symbols. The image should convey the fun and excitement of the game and appeal to players who enjoy playful and adventurous slot games.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Discover the amusing online slot game Ace Ventura and its various bonus features. Play for free and enjoy the immersive experience.", 2)

Write-Output "done"
